$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "62.765.78"
Set-TextValue $ws.Range("E2") "  -0.67%  "

Set-TextValue $ws.Range("D3") "2.574.15"
Set-TextValue $ws.Range("E3") "  +0.60%  "

Set-TextValue $ws.Range("E4") "  -0.05%  "

Set-TextValue $ws.Range("D5") "581.23"
Set-TextValue $ws.Range("E5") "  -0.44%  "

Set-TextValue $ws.Range("D6") "143.65"
Set-TextValue $ws.Range("E6") "  -2.92%  "

Set-TextValue $ws.Range("E7") "  -0.02%  "

Set-TextValue $ws.Range("D8") "0.589"
Set-TextValue $ws.Range("E8") "  +0.53%  "

Set-TextValue $ws.Range("E9") "  -2.49%  "

Set-TextValue $ws.Range("D10") "5.59"
Set-TextValue $ws.Range("E10") "  -0.47%  "

Set-TextValue $ws.Range("E11") "  -0.52%  "

Set-TextValue $ws.Range("E12") "  -1.82%  "

Set-TextValue $ws.Range("D13") "26.84"
Set-TextValue $ws.Range("E13") "  -2.83%  "

Set-TextValue $ws.Range("D14") "3.032.65"
Set-TextValue $ws.Range("E14") "  +0.55%  "

Set-TextValue $ws.Range("D15") "62.709.85"
Set-TextValue $ws.Range("E15") "  -0.55%  "

Set-TextValue $ws.Range("E16") "  -2.62%  "

Set-TextValue $ws.Range("D17") "2.576.50"
Set-TextValue $ws.Range("E17") "  +0.29%  "

Set-TextValue $ws.Range("E18") "  -2.78%  "

Set-TextValue $ws.Range("D19") "340.63"
Set-TextValue $ws.Range("E19") "  -0.34%  "

Set-TextValue $ws.Range("E20") "  -2.14%  "

Set-TextValue $ws.Range("E21") "  -2.53%  "

Set-TextValue $ws.Range("D22") "1.00"
Set-TextValue $ws.Range("E22") "  +0.07%  "

Set-TextValue $ws.Range("D23") "66.83"
Set-TextValue $ws.Range("E23") "  +0.47%  "

Set-TextValue $ws.Range("D24") "1.58"
Set-TextValue $ws.Range("E24") "  -4.13%  "

Set-TextValue $ws.Range("B25") "SuiNetwork"
Set-TextValue $ws.Range("C25") "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D25") "1.52"
Set-TextValue $ws.Range("E25") "  +1.94%  "

Set-TextValue $ws.Range("B26") "Kaspa"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D26") "0.164"
Set-TextValue $ws.Range("E26") "  -3.89%  "

Set-TextValue $ws.Range("D27") "1.00"
Set-TextValue $ws.Range("E27") "  +0.01%  "

Set-TextValue $ws.Range("D28") "7.89"
Set-TextValue $ws.Range("E28") "  -2.97%  "

Set-TextValue $ws.Range("D29") "8.21"
Set-TextValue $ws.Range("E29") "  -4.04%  "

Set-TextValue $ws.Range("D30") "1.93"
Set-TextValue $ws.Range("E30") "  -3.49%  "

Set-TextValue $ws.Range("D31") "454.17"
Set-TextValue $ws.Range("E31") "  +2.77%  "

Set-TextValue $ws.Range("E32") "  -3.64%  "

Set-TextValue $ws.Range("E33") "  +1.17%  "

Set-TextValue $ws.Range("D34") "176.33"
Set-TextValue $ws.Range("E34") "  -0.34%  "

Set-TextValue $ws.Range("E35") "  +0.15%  "

Set-TextValue $ws.Range("E36") "  -2.23%  "

Set-TextValue $ws.Range("E37") "  -2.28%  "

Set-TextValue $ws.Range("E38") "  -1.51%  "

Set-TextValue $ws.Range("E39") "  +0.01%  "

Set-TextValue $ws.Range("E40") "  -3.28%  "

Set-TextValue $ws.Range("D41") "40.02"
Set-TextValue $ws.Range("E41") "  +0.87%  "

Set-TextValue $ws.Range("D42") "156.98"
Set-TextValue $ws.Range("E42") "  +4.10%  "

Set-TextValue $ws.Range("D43") "3.69"
Set-TextValue $ws.Range("E43") "  -3.92%  "

Set-TextValue $ws.Range("D44") "0.632"
Set-TextValue $ws.Range("E44") "  +3.35%  "

Set-TextValue $ws.Range("D45") "21.04"
Set-TextValue $ws.Range("E45") "  -0.30%  "

Set-TextValue $ws.Range("E46") "  -3.32%  "

Set-TextValue $ws.Range("D47") "0.0957"
Set-TextValue $ws.Range("E47") "  -1.91%  "

Set-TextValue $ws.Range("E48") "  -3.14%  "

Set-TextValue $ws.Range("D49") "17.93"
Set-TextValue $ws.Range("E49") "  -2.65%  "

Set-TextValue $ws.Range("E51") "  -4.11%  "
